$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename sheets: Sheet1 -> work_schedule, Sheet2 -> r_users
# ---------------------------------------------------------------------------
$wsSchedule = $wb.Worksheets.Item(1)
$wsSchedule.Name = "work_schedule"

$wsUsers = $wb.Worksheets.Item(2)
$wsUsers.Name = "r_users"

# ---------------------------------------------------------------------------
# 2. work_schedule (Sheet1): row 11 content changed from the "Intro to
#    Machine Learning" webinar to a new "GTFS data in R" webinar.
# ---------------------------------------------------------------------------
$wsSchedule.Range("E11").Value = "GTFS data in R"
$wsSchedule.Range("H11").Value = "topics/gtfs_data.html"

# ---------------------------------------------------------------------------
# 3. r_users (Sheet2): build out the new roster table.
#    Cells are written in the order that reproduces the target shared-string
#    layout (header row first, then first/last name pairs per row, with
#    r_version written after the first five data rows).
# ---------------------------------------------------------------------------
$wsUsers.Range("A1").Value = "user_id"
$wsUsers.Range("B1").Value = "first_name"
$wsUsers.Range("B2").Value = "Bryan"
$wsUsers.Range("C1").Value = "last_name"
$wsUsers.Range("C2").Value = "Blanc"
$wsUsers.Range("B3").Value = "Esther"
$wsUsers.Range("C3").Value = "Needham"
$wsUsers.Range("B4").Value = "Oren"
$wsUsers.Range("C4").Value = "Eshel"
$wsUsers.Range("B5").Value = "Paul"
$wsUsers.Range("C5").Value = "Leitman"
$wsUsers.Range("D1").Value = "r_version"
$wsUsers.Range("B6").Value = "Joseph"
$wsUsers.Range("C6").Value = "Poirier"
$wsUsers.Range("B7").Value = "Tomoko"
$wsUsers.Range("C7").Value = "DeLaTorre"

$wsUsers.Range("A2").Value = 1
$wsUsers.Range("A3").Value = 2
$wsUsers.Range("A4").Value = 3
$wsUsers.Range("A5").Value = 4
$wsUsers.Range("A6").Value = 5
$wsUsers.Range("A7").Value = 6
$wsUsers.Range("A8").Value = 7

# Bold header row (no border) - matches the new cellXfs style used on row 1.
$wsUsers.Range("A1:D1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 4. Selections / active sheet bookkeeping.
# ---------------------------------------------------------------------------
$wsSchedule.Range("H12").Select()
$wsUsers.Activate()
$wsUsers.Range("B8").Select()
